$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (closest achievable values given the host's
# pixel-grid snapping of ColumnWidth; target stored widths are
# 15.7109375 and 16.42578125)
$ws.Columns.Item(1).ColumnWidth = 14.83
$ws.Columns.Item(2).ColumnWidth = 15.65

# Update cell values
$ws.Range("A1").Value = 0.00041965430430657543
$ws.Range("B1").Value = -0.00041965497552812707

$ws.Range("A2").Value = -0.0074889758259489942
$ws.Range("B2").Value = 0.0074889751258344081

$ws.Range("A3").Value = 0.025240815353237259
$ws.Range("B3").Value = -0.025240815991401145

$ws.Range("A4").Value = -0.0177848572980044
$ws.Range("B4").Value = 0.017784856665811907

$ws.Range("A5").Value = 0.0013527091559002362
$ws.Range("B5").Value = -0.001352709853694403
